$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column (value + formatting copied from an existing header cell)
$ws.Range("H1").Value = "Recomendação Investimento"
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Values copied from the diff (column H, rows 2-31)
$values = @(
    1024.199981689453,
    188.0999994277954,
    356.1199922561646,
    160.950005531311,
    1142.309982299805,
    1069.809982299805,
    1664.600044250488,
    642.9300022125244,
    116.7000031471252,
    912.0500221252441,
    1364.449977874756,
    333.7900066375732,
    282.1699867248535,
    514.75,
    57.4200005531311,
    545.1999778747559,
    759.5100154876709,
    111.9399969577789,
    685.8499889373779,
    738.0500221252441,
    223.299994468689,
    329.1500110626221,
    696,
    51.61999917030334,
    327.1199922561646,
    350.9000110626221,
    116,
    198.0699977874756,
    305.950005531311,
    230.549994468689
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
